# updated rate data from CDC
#
# CDC refreshed its weekly COVID-19 case-rate-by-age-group figures, which
# revises the last 7 already-published weeks (rows 76-82 on "Sheet1") and
# adds the newly published week ending 2021-09-25 (row 83). "Sheet2" keeps a
# mirror of just those 8 weeks, so it is rebuilt from the same data.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Columns: week, a0_4, a5_11, a12_15, a16_17, a18_29, a30_39, a40_49, a50_64, a65_74, a75
$data = @(
  @(44415, 96.1, 138,   155.3, 189.2, 234.9, 228.8, 196.4, 146.3, 108.3, 103.7),
  @(44422, 122.5,192.1, 207.3, 241.2, 255.5, 255.6, 221.1, 162.7, 121.7, 110.6),
  @(44429, 131.6,233.5, 257.5, 269.7, 241.8, 248.2, 212.4, 158.6, 123.4, 110.7),
  @(44436, 162.5,292.8, 337.9, 349.3, 294.8, 298.3, 256.2, 191.3, 145.6, 130.1),
  @(44443, 136,  259.8, 301.2, 312.1, 242.3, 249.1, 215.9, 161.3, 125.2, 113.8),
  @(44450, 131.6,246.2, 289.2, 302.2, 240.7, 247.9, 217.3, 162.4, 125.3, 114.3),
  @(44457, 150,  265.2, 289.2, 307.1, 268.2, 281.7, 250,   193.1, 149.7, 137.4),
  @(44464, 111.8,202.9, 211,   219.5, 181.7, 199.4, 177.8, 138.1, 112.5, 104.9)
)

# Sheet1: overwrite the existing rows 76-82 and add the new row 83.
$firstRow = 76
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $firstRow + $i
    $values = $data[$i]
    for ($c = 1; $c -le $values.Length; $c++) {
        $ws1.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}

# Sheet2: same 8 weeks of data, written starting at row 1.
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $values = $data[$i]
    for ($c = 1; $c -le $values.Length; $c++) {
        $ws2.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}

# Restore each sheet's on-screen selection/scroll position.
$ws1.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 53
$ws1.Range("C87").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("A1:K8").Select() | Out-Null

$ws1.Activate() | Out-Null
